$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2376.037
$ws.Range("J17").Value = 2376.037
$ws.Range("L17").Value = 7128.110999999999
$ws.Range("N17").Value = -7464.110999999999
$ws.Range("H70").Value = 1303.88
$ws.Range("J70").Value = 1273.1333
$ws.Range("L70").Value = 3819.3999
$ws.Range("N70").Value = -4359.3999
$ws.Range("H73").Value = 1303.88
$ws.Range("J73").Value = 1273.1333
$ws.Range("L73").Value = 3819.3999
$ws.Range("N73").Value = -5691.3999
$ws.Range("H106").Value = 72606.91
$ws.Range("I106").Value = 1988.7778
$ws.Range("K106").Value = 1988.7778
$ws.Range("M106").Value = -1357.7778
$ws.Range("H112").Value = 2332.0952
$ws.Range("I112").Value = 824.6667
$ws.Range("J112").Value = 2583.3333
$ws.Range("K112").Value = 2474.0001
$ws.Range("L112").Value = 7749.999899999999
$ws.Range("M112").Value = -1366.0001
$ws.Range("N112").Value = -9965.999899999999
$ws.Range("H129").Value = 886.7
$ws.Range("I129").Value = 546
$ws.Range("J129").Value = 920.3955999999999
$ws.Range("K129").Value = 1638
$ws.Range("L129").Value = 2761.1868
$ws.Range("M129").Value = 3362
$ws.Range("N129").Value = -12761.1868
$ws.Range("H132").Value = 47729.906
$ws.Range("I132").Value = 6859
$ws.Range("J132").Value = 149907.17
$ws.Range("K132").Value = 20577
$ws.Range("L132").Value = 449721.51
$ws.Range("M132").Value = -18047
$ws.Range("N132").Value = -454781.51
$ws.Range("H138").Value = 1952.5977
$ws.Range("I138").Value = 1437.9131
$ws.Range("J138").Value = 2137.5625
$ws.Range("K138").Value = 4313.7393
$ws.Range("L138").Value = 6412.6875
$ws.Range("M138").Value = 826.2606999999998
$ws.Range("N138").Value = -16692.6875

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 1000
$ws.Range("K22").Value = 1000
$ws.Range("M22").Value = -701
$ws.Range("H32").Value = 10527.712
$ws.Range("I32").Value = 9675.404
$ws.Range("K32").Value = 9675.404
$ws.Range("M32").Value = -9388.404
$ws.Range("H45").Value = 2110.1155
$ws.Range("I45").Value = 1957.7273
$ws.Range("K45").Value = 1957.7273
$ws.Range("M45").Value = -1580.7273
$ws.Range("H61").Value = 2006.0834
$ws.Range("I61").Value = 1526.0588
$ws.Range("K61").Value = 1526.0588
$ws.Range("M61").Value = -1314.0588
$ws.Range("H74").Value = 1988.659
$ws.Range("I74").Value = 1677.8611
$ws.Range("J74").Value = 3387.25
$ws.Range("K74").Value = 1677.8611
$ws.Range("L74").Value = 3387.25
$ws.Range("M74").Value = -803.8611000000001
$ws.Range("N74").Value = -5135.25
$ws.Range("H77").Value = 1988.659
$ws.Range("I77").Value = 1677.8611
$ws.Range("J77").Value = 3387.25
$ws.Range("K77").Value = 8389.3055
$ws.Range("L77").Value = 16936.25
$ws.Range("M77").Value = -4021.3055
$ws.Range("N77").Value = -25672.25
$ws.Range("H88").Value = 11063794
$ws.Range("I88").Value = 25003400
$ws.Range("J88").Value = 2485575.8
$ws.Range("K88").Value = 25003400
$ws.Range("L88").Value = 2485575.8
$ws.Range("M88").Value = -25002994
$ws.Range("N88").Value = -2486387.8
$ws.Range("H91").Value = 11063794
$ws.Range("I91").Value = 25003400
$ws.Range("J91").Value = 2485575.8
$ws.Range("K91").Value = 25003400
$ws.Range("L91").Value = 2485575.8
$ws.Range("M91").Value = -25001996
$ws.Range("N91").Value = -2488383.8
$ws.Range("H102").Value = 14287
$ws.Range("I102").Value = 2201.4583
$ws.Range("J102").Value = 38458.082
$ws.Range("K102").Value = 2201.4583
$ws.Range("L102").Value = 38458.082
$ws.Range("M102").Value = -579.4582999999998
$ws.Range("N102").Value = -41702.082
$ws.Range("H110").Value = 1263.762
$ws.Range("I110").Value = 1180.7028
$ws.Range("J110").Value = 1878.4
$ws.Range("K110").Value = 1180.7028
$ws.Range("L110").Value = 1878.4
$ws.Range("M110").Value = 864.2972
$ws.Range("N110").Value = -5968.4
$ws.Range("H132").Value = 16668254
$ws.Range("I132").Value = 20000944
$ws.Range("J132").Value = 4799.2
$ws.Range("K132").Value = 60002832
$ws.Range("L132").Value = 14397.6
$ws.Range("M132").Value = -60000302
$ws.Range("N132").Value = -19457.6
$ws.Range("H136").Value = 2006.0834
$ws.Range("I136").Value = 1526.0588
$ws.Range("K136").Value = 4578.1764
$ws.Range("M136").Value = -2028.1764

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1859.5
$ws.Range("I86").Value = 2161.2
$ws.Range("J86").Value = 1557.8
$ws.Range("K86").Value = 2161.2
$ws.Range("L86").Value = 1557.8
$ws.Range("M86").Value = -1038.2
$ws.Range("N86").Value = -3803.8
$ws.Range("H89").Value = 1859.5
$ws.Range("I89").Value = 2161.2
$ws.Range("J89").Value = 1557.8
$ws.Range("K89").Value = 10806
$ws.Range("L89").Value = 7789
$ws.Range("M89").Value = -5190
$ws.Range("N89").Value = -19021
$ws.Range("H105").Value = 3477.8462
$ws.Range("I105").Value = 2108
$ws.Range("K105").Value = 2108
$ws.Range("M105").Value = -361

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("M56").ClearContents()
$ws.Range("N56").ClearContents()
$ws.Range("H122").Value = 120649.4
$ws.Range("I122").Value = 133943.78
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 401831.34
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -399381.34
$ws.Range("N122").Value = -7900
$ws.Range("H132").Value = 111470.62
$ws.Range("I132").Value = 1401.5
$ws.Range("J132").Value = 205815.58
$ws.Range("K132").Value = 4204.5
$ws.Range("L132").Value = 617446.74
$ws.Range("M132").Value = -1674.5
$ws.Range("N132").Value = -622506.74
$ws.Range("H134").Value = 1079174.8
$ws.Range("I134").Value = 1145.125
$ws.Range("J134").Value = 2804022.2
$ws.Range("K134").Value = 3435.375
$ws.Range("L134").Value = 8412066.600000001
$ws.Range("M134").Value = -900.375
$ws.Range("N134").Value = -8417136.600000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5471.8125
$ws.Range("I3").Value = 1778.5714
$ws.Range("J3").Value = 8344.333000000001
$ws.Range("K3").Value = 5335.7142
$ws.Range("L3").Value = 25032.999
$ws.Range("M3").Value = -5223.7142
$ws.Range("N3").Value = -25256.999
$ws.Range("H104").Value = 1492
$ws.Range("J104").Value = 1484
$ws.Range("L104").Value = 4452
$ws.Range("N104").Value = -9694
$ws.Range("H113").Value = 6705.5293
$ws.Range("I113").Value = 8885.833000000001
$ws.Range("J113").Value = 1472.8
$ws.Range("K113").Value = 26657.499
$ws.Range("L113").Value = 4418.4
$ws.Range("M113").Value = -24487.499
$ws.Range("N113").Value = -8758.4
$ws.Range("H121").Value = 92116.32000000001
$ws.Range("I121").Value = 610
$ws.Range("J121").Value = 100971.77
$ws.Range("K121").Value = 1830
$ws.Range("L121").Value = 302915.31
$ws.Range("M121").Value = -520
$ws.Range("N121").Value = -305535.31

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 45203.332
$ws.Range("J53").Value = 45203.332
$ws.Range("L53").Value = 45203.332
$ws.Range("N53").Value = -46465.332
$ws.Range("H126").Value = 8218.944
$ws.Range("I126").Value = 13742.111
$ws.Range("J126").Value = 2695.7778
$ws.Range("K126").Value = 41226.333
$ws.Range("L126").Value = 8087.3334
$ws.Range("M126").Value = -38756.333
$ws.Range("N126").Value = -13027.3334

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2247.0527
$ws.Range("I122").Value = 2294
$ws.Range("J122").Value = 1996.6666
$ws.Range("K122").Value = 6882
$ws.Range("L122").Value = 5989.9998
$ws.Range("M122").Value = -4432
$ws.Range("N122").Value = -10889.9998
$ws.Range("H132").Value = 2727.9473
$ws.Range("I132").Value = 1933.6957
$ws.Range("J132").Value = 3945.8
$ws.Range("K132").Value = 5801.0871
$ws.Range("L132").Value = 11837.4
$ws.Range("M132").Value = -3271.0871
$ws.Range("N132").Value = -16897.4
$ws.Range("H136").Value = 2465.6667
$ws.Range("I136").Value = 1698.8667
$ws.Range("J136").Value = 4382.6665
$ws.Range("K136").Value = 5096.6001
$ws.Range("L136").Value = 13147.9995
$ws.Range("M136").Value = -2546.6001
$ws.Range("N136").Value = -18247.9995

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("H107").Value = 7143714.5
$ws.Range("I107").Value = 320.3
$ws.Range("J107").Value = 25002200
$ws.Range("K107").Value = 960.9000000000001
$ws.Range("L107").Value = 75006600
$ws.Range("M107").Value = 959.0999999999999
$ws.Range("N107").Value = -75010440
$ws.Range("H113").Value = 410.3
$ws.Range("I113").Value = 410.3
$ws.Range("K113").Value = 1230.9
$ws.Range("M113").Value = 939.0999999999999
$ws.Range("H122").Value = 7143807
$ws.Range("I122").Value = 14286164
$ws.Range("J122").Value = 1450
$ws.Range("K122").Value = 42858492
$ws.Range("L122").Value = 4350
$ws.Range("M122").Value = -42856042
$ws.Range("N122").Value = -9250
$ws.Range("H136").Value = 371862.03
$ws.Range("I136").Value = 626073.6
$ws.Range("J136").Value = 2099.7273
$ws.Range("K136").Value = 1878220.8
$ws.Range("L136").Value = 6299.1819
$ws.Range("M136").Value = -1875670.8
$ws.Range("N136").Value = -11399.1819
